# Add "Beverly's Output" column (H) + a blank "Status" column (I) to the
# Temperature Converter test-case sheet, matching Beverly's manual run of
# the test suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("H1").Value = "Beverly's Output"
$ws.Range("I1").Value = "Status"

# --- Beverly's Output values, row by row ------------------------------
# (numeric results typed as numbers; textual results as strings)
$ws.Range("H2").Value = 273.14999999999998
$ws.Range("H3").Value = 373.15

# Rows 5 and 6 were typed before row 4, which is why their strings land
# earlier in the shared-string table than "0=".
$ws.Range("H5").Value = "100="
$ws.Range("H6").Value = "273.15="
$ws.Range("H4").Value = "0="

$ws.Range("H7").Value = 373.15
$ws.Range("H8").Value = 32
$ws.Range("H9").Value = -459.67

$ws.Range("H13").Value = "invalid input"

$ws.Range("H14").Value = 273.14999999999998
$ws.Range("H15").Value = 273.14999999999998
$ws.Range("H16").Value = 300
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 32
$ws.Range("H19").Value = 300

# --- New "Status" (I) column: every data row gets the bright-green fill
# used to flag Beverly's results, left blank for manual marking. --------
$greenFill = 5296274   # RGB(146, 208, 80) == FF92D050

for ($r = 2; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Interior.Color = $greenFill
}

# --- Selection / scroll state matching the saved workbook --------------
$ws.Range("H4").Select()
